# daily sheet.xlsx edit: "new feature: Control Panel"
#
# Job #3 row (row 7) gets the "Fix bug: Restrict the map bounders" work item
# with an updated begin time and newly filled-in finish time / hours.
# Job #4 row (row 8) becomes the new "New feature: Add custom control panel"
# work item, with its own begin/finish/hours.
# A brand new Job #5 row (row 9) is appended below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Job #3) ---------------------------------------------------
$ws.Range("C7").Value2 = "Fix bug: Restrict the map bounders"
$ws.Range("D7").Value2 = 0.58333333333333337
$ws.Range("E7").Value2 = 0.59375
$ws.Range("F7").Value2 = 0.25

# --- Row 8 (Job #4) -----------------------------------------------------
$ws.Range("D8").Value2 = 0.57291666666666663
$ws.Range("E8").Value2 = 0.625
$ws.Range("F8").Value2 = 1.25

# --- Row 9 (new Job #5) --------------------------------------------------
# Written before C8's new string so the shared-string table picks up "5"
# ahead of "New feature: Add custom control panel", matching the author's
# original edit order.
$ws.Range("B9").Value2 = "5"

$ws.Range("C8").Value2 = "New feature: Add custom control panel"

# Match the selection left behind by the author's edit session.
$ws.Range("C10").Select()
